$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts two new price-report rows (new dates 2021-10-xx
# range) ahead of the existing row that used to be row 14, pushing all
# subsequent rows down by two positions in total.
$ws.Rows("14:14").Insert()
$ws.Rows("16:16").Insert()

# New row 14 (Comercializadora del Agro de Limarí - Arveja Verde)
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44482
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112022
$ws.Range("G14").Value = "Arveja Verde"
$ws.Range("H14").Value = "Perfection"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 18000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 19000
$ws.Range("N14").Value = "$/malla 25 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 760
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"

# New row 16 (Comercializadora del Agro de Limarí - Arveja Verde)
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44483
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 100112022
$ws.Range("G16").Value = "Arveja Verde"
$ws.Range("H16").Value = "Perfection"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 19000
$ws.Range("N16").Value = "$/malla 25 kilos"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 760
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"
